$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (2007年 data), shifting all rows below it up by one.
$ws.Rows.Item(2).Delete()
